# Update the NATMI ligand-receptor pair sheet (Angpt1-Itga5) with the
# recomputed TPM-based values. A new sending/target cluster, "Resolving-Mac",
# is added, expanding the 3x4 cluster table to a full 4x4 cross of
# {ECs, FAPs, MuSCs, Resolving-Mac}, and every metric column is refreshed
# with the new TPM-derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('ECs', 'Angpt1', 'Itga5', 'ECs', 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 35.04689966666667, 105.140699, 0.3824629895491901, 0.3824629895491901, 6.261771151943889, 56.355940367495, 0.002180395164793063, 0.002180395164793062),
  @('ECs', 'Angpt1', 'Itga5', 'FAPs', 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 29.913269, 89.739807, 0.3264402385872224, 0.3264402385872223, 5.344553916781666, 48.100985251035, 0.001861013319611491, 0.001861013319611491),
  @('ECs', 'Angpt1', 'Itga5', 'MuSCs', 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 8.911727666666666, 26.735183, 0.09725271102035077, 0.09725271102035075, 1.592243529323889, 14.330191763915, 0.0005544310081394615, 0.0005544310081394614),
  @('ECs', 'Angpt1', 'Itga5', 'Resolving-Mac', 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 17.76285166666667, 53.288555, 0.1938440608432367, 0.1938440608432367, 3.173659102530555, 28.562931922775, 0.001105091641637356, 0.001105091641637356),
  @('FAPs', 'Angpt1', 'Itga5', 'ECs', 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 35.04689966666667, 105.140699, 0.3824629895491901, 0.3824629895491901, 776.2128225735896, 6985.915403162307, 0.2702830627504536, 0.2702830627504536),
  @('FAPs', 'Angpt1', 'Itga5', 'FAPs', 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 29.913269, 89.739807, 0.3264402385872224, 0.3264402385872223, 662.5140364406286, 5962.626327965658, 0.2306923020037615, 0.2306923020037615),
  @('FAPs', 'Angpt1', 'Itga5', 'MuSCs', 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 8.911727666666666, 26.735183, 0.09725271102035077, 0.09725271102035075, 197.3754412499335, 1776.378971249402, 0.06872759277008286, 0.06872759277008285),
  @('FAPs', 'Angpt1', 'Itga5', 'Resolving-Mac', 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 17.76285166666667, 53.288555, 0.1938440608432367, 0.1938440608432367, 393.4086427123522, 3540.67778441117, 0.1369878076894467, 0.1369878076894466),
  @('MuSCs', 'Angpt1', 'Itga5', 'ECs', 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 35.04689966666667, 105.140699, 0.3824629895491901, 0.3824629895491901, 314.2769979648911, 2828.49298168402, 0.1094335819915107, 0.1094335819915107),
  @('MuSCs', 'Angpt1', 'Itga5', 'FAPs', 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 29.913269, 89.739807, 0.3264402385872224, 0.3264402385872223, 268.2420547908733, 2414.17849311786, 0.09340387329208115, 0.09340387329208112),
  @('MuSCs', 'Angpt1', 'Itga5', 'MuSCs', 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 8.911727666666666, 26.735183, 0.09725271102035077, 0.09725271102035075, 79.91437315137111, 719.22935836234, 0.02782677753444023, 0.02782677753444022),
  @('MuSCs', 'Angpt1', 'Itga5', 'Resolving-Mac', 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 17.76285166666667, 53.288555, 0.1938440608432367, 0.1938440608432367, 159.2852934265444, 1433.5676408389, 0.05546432074606644, 0.05546432074606641),
  @('Resolving-Mac', 'Angpt1', 'Itga5', 'ECs', 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 35.04689966666667, 105.140699, 0.3824629895491901, 0.3824629895491901, 1.625323336641445, 14.627910029773, 0.0005659496424327469, 0.0005659496424327467),
  @('Resolving-Mac', 'Angpt1', 'Itga5', 'FAPs', 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 29.913269, 89.739807, 0.3264402385872224, 0.3264402385872223, 1.387247792054333, 12.485230128489, 0.0004830499717681515, 0.0004830499717681513),
  @('Resolving-Mac', 'Angpt1', 'Itga5', 'MuSCs', 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 8.911727666666666, 26.735183, 0.09725271102035077, 0.09725271102035075, 0.4132873116934444, 3.719585805241, 0.0001439097076882098, 0.0001439097076882097),
  @('Resolving-Mac', 'Angpt1', 'Itga5', 'Resolving-Mac', 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 17.76285166666667, 53.288555, 0.1938440608432367, 0.1938440608432367, 0.8237640879427778, 7.413876791485, 0.0002868407660862875, 0.0002868407660862874),
)

# Row 1 (headers) is untouched. Data rows start at row 2; write all 20
# columns (A:T) for each of the 16 rows (4 sending clusters x 4 target
# clusters), growing the sheet from A1:T13 to A1:T17.
$rowCount = $data.Count
$colCount = 20

for ($i = 0; $i -lt $rowCount; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colCount; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
